# "Generate Report for Handoff"
# Updates the localization-status report after a handoff was generated:
#   - Status moves from "In Translation" to "Ready for handoff"
#   - The handoff timestamps are refreshed
#   - The Status columns are widened slightly (to fit the longer text)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status column

# --- Refreshed handoff datetimes --------------------------------------------
$wsOverview.Range("G2").Value = "2016-08-17 06:52:03"  # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value = "2016-08-17 06:52:03"      # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value = "2016-08-17 06:51:56"      # Latest Handoff Datetime (zh-cn)

# --- Widen the Status columns to fit the new, longer status text ------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # column C (Status)
